# Add a new "2022-Q4" sheet (with its fund holdings data) positioned right
# after "总计" and before "2022-Q3", and update the "总计" (totals) sheet so
# it gains a new leading row summarising the 2022-Q4 quarter.

$wb = $excel.ActiveWorkbook

# --- locate existing sheets by position (stable only until the next
#     structural change - worksheet handles are positional in this COM
#     model, so we always re-resolve sheets we need *after* an insert by
#     name rather than by a previously-captured variable). -------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3Pos = $wb.Worksheets.Item(2)

# --- insert the new quarter sheet right before "2022-Q3" ---------------
$wsQ4 = $wb.Worksheets.Add($wsQ3Pos)
$wsQ4.Name = "2022-Q4"

# re-resolve the Q3 sheet by name now that the sheet collection changed
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# --- copy header row + index column formatting from the Q3 sheet -------
$wsQ3.Range("A1:H1").Copy()
$wsQ4.Range("A1:H1").PasteSpecial(-4122)

$wsQ3.Range("A2:A8").Copy()
$wsQ4.Range("A2:A8").PasteSpecial(-4122)

# --- header labels -------------------------------------------------------
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# --- index column (A2:A8) ------------------------------------------------
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("A4").Value = 2
$wsQ4.Range("A5").Value = 3
$wsQ4.Range("A6").Value = 4
$wsQ4.Range("A7").Value = 5
$wsQ4.Range("A8").Value = 6

# --- fund data - B,C,D,E,F,G are stored as text ---------------------------
$wsQ4.Range("B2:G8").NumberFormat = "@"

$wsQ4.Range("B2").Value = "202027"
$wsQ4.Range("C2").Value = "南方高端装备灵活配置混合A"
$wsQ4.Range("D2").Value = "16.75"
$wsQ4.Range("E2").Value = "93.09"
$wsQ4.Range("F2").Value = "7.90"
$wsQ4.Range("G2").Value = "1.3232"
$wsQ4.Range("H2").Value = 3

$wsQ4.Range("B3").Value = "005207"
$wsQ4.Range("C3").Value = "南方高端装备灵活配置混合C"
$wsQ4.Range("D3").Value = "5.35"
$wsQ4.Range("E3").Value = "93.09"
$wsQ4.Range("F3").Value = "7.90"
$wsQ4.Range("G3").Value = "0.4226"
$wsQ4.Range("H3").Value = 3

$wsQ4.Range("B4").Value = "017732"
$wsQ4.Range("C4").Value = "鹏华核心优势混合C"
$wsQ4.Range("D4").Value = "5.55"
$wsQ4.Range("E4").Value = "88.10"
$wsQ4.Range("F4").Value = "3.49"
$wsQ4.Range("G4").Value = "0.1937"
$wsQ4.Range("H4").Value = 7

$wsQ4.Range("B5").Value = "001223"
$wsQ4.Range("C5").Value = "鹏华文化传媒娱乐股票"
$wsQ4.Range("D5").Value = "0.77"
$wsQ4.Range("E5").Value = "83.57"
$wsQ4.Range("F5").Value = "4.26"
$wsQ4.Range("G5").Value = "0.0328"
$wsQ4.Range("H5").Value = 6

$wsQ4.Range("B6").Value = "016013"
$wsQ4.Range("C6").Value = "南方碳中和股票A"
$wsQ4.Range("D6").Value = "0.42"
$wsQ4.Range("E6").Value = "84.75"
$wsQ4.Range("F6").Value = "5.72"
$wsQ4.Range("G6").Value = "0.0240"
$wsQ4.Range("H6").Value = 2

$wsQ4.Range("B7").Value = "016014"
$wsQ4.Range("C7").Value = "南方碳中和股票C"
$wsQ4.Range("D7").Value = "0.01"
$wsQ4.Range("E7").Value = "84.75"
$wsQ4.Range("F7").Value = "5.72"
$wsQ4.Range("G7").Value = "0.0006"
$wsQ4.Range("H7").Value = 2

$wsQ4.Range("B8").Value = "006976"
$wsQ4.Range("C8").Value = "鹏华核心优势混合A"
$wsQ4.Range("D8").Value = "0.00"
$wsQ4.Range("E8").Value = "88.10"
$wsQ4.Range("F8").Value = "3.49"
# G8 holds a real 0 number (not text) per the source data
$wsQ4.Range("G8").ClearFormats()
$wsQ4.Range("G8").Value = 0
$wsQ4.Range("H8").Value = 7

# --- update the "总计" (totals) sheet: insert the Q4 totals as the new
#     first data row, and push the Q3 / Q1 rows down -----------------
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3:A4").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.09

$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 15
$wsTotal.Range("D3").Value = 3.19

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 2
